$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "29.20", "6.80") would be
# auto-coerced to a numeric value by Excel (losing the literal formatting, e.g.
# trailing zeros) unless the cell is explicitly forced to Text first. The source
# sheet stores these as plain inline strings with the default (General) style, so
# after writing the text we flip the style back to "Normal" to avoid leaving a
# stray Text number-format on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "65.971.23"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.667.10"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "599.79"
$ws.Range("E5").Value = "  -0.15%  "
Set-TextValue $ws.Range("D6") "159.06"
$ws.Range("E6").Value = "  +1.43%  "
Set-TextValue $ws.Range("D7") "0.653"
$ws.Range("E7").Value = "  +4.90%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("E11").Value = "  -0.08%  "
Set-TextValue $ws.Range("D12") "0.158"
$ws.Range("E12").Value = "  +1.71%  "
Set-TextValue $ws.Range("D13") "29.20"
$ws.Range("E13").Value = "  -0.63%  "
Set-TextValue $ws.Range("D14") "0.0000197"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "3.145.83"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "65.807.26"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").Value = "2.657.64"
$ws.Range("E17").Value = "  -0.56%  "
Set-TextValue $ws.Range("D18") "12.67"
$ws.Range("E18").Value = "  -1.90%  "
Set-TextValue $ws.Range("D19") "4.83"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D20") "355.02"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D21") "7.53"
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("E22").Value = "  -0.06%  "
Set-TextValue $ws.Range("D23") "70.14"
$ws.Range("E23").Value = "  +0.51%  "
Set-TextValue $ws.Range("D24") "1.83"
$ws.Range("E24").Value = "  +11.15%  "
Set-TextValue $ws.Range("D25") "0.0000114"
$ws.Range("E25").Value = "  +1.76%  "
Set-TextValue $ws.Range("D26") "9.72"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("E27").Value = "  +2.15%  "
Set-TextValue $ws.Range("D28") "584.24"
$ws.Range("E28").Value = "  +10.68%  "
Set-TextValue $ws.Range("D29") "8.22"
$ws.Range("E29").Value = "  +1.90%  "
Set-TextValue $ws.Range("D30") "0.164"
$ws.Range("E30").Value = "  -2.21%  "
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("E33").Value = "  +3.00%  "
Set-TextValue $ws.Range("D34") "6.80"
$ws.Range("E34").Value = "  +4.69%  "
Set-TextValue $ws.Range("D35") "5.57"
$ws.Range("E35").Value = "  +1.25%  "
Set-TextValue $ws.Range("D36") "0.425"
$ws.Range("E36").Value = "  +0.11%  "
Set-TextValue $ws.Range("D37") "20.66"
$ws.Range("E37").Value = "  +0.25%  "
Set-TextValue $ws.Range("D38") "1.98"
$ws.Range("E38").Value = "  +1.87%  "
$ws.Range("E39").Value = "  -0.07%  "
Set-TextValue $ws.Range("D40") "154.57"
$ws.Range("E40").Value = "  -2.10%  "
Set-TextValue $ws.Range("D41") "162.47"
$ws.Range("E41").Value = "  -1.09%  "
Set-TextValue $ws.Range("D42") "4.13"
$ws.Range("E42").Value = "  -0.44%  "
Set-TextValue $ws.Range("D43") "2.40"
$ws.Range("E43").Value = "  +3.16%  "
Set-TextValue $ws.Range("D44") "0.0622"
$ws.Range("E44").Value = "  +1.42%  "
Set-TextValue $ws.Range("D45") "23.36"
$ws.Range("E45").Value = "  +2.08%  "
Set-TextValue $ws.Range("D46") "0.647"
$ws.Range("E46").Value = "  +0.81%  "
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("E48").Value = "  +1.90%  "
Set-TextValue $ws.Range("D49") "19.93"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").Value = "0.0₆0248"
$ws.Range("E50").Value = "  -7.41%  "
Set-TextValue $ws.Range("D51") "0.824"
$ws.Range("E51").Value = "  +0.79%  "
